$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Narrow column A slightly (closest value reachable through Excel's
# pixel-grid rounding of ColumnWidth) to match the new width of 13.7109375
$ws.Columns.Item(1).ColumnWidth = 12.8333333333333

# Update the six data values in rows 1-3 (columns A and B); rows 4-5 stay untouched
$ws.Range("A1").Value = 0.050725434207982312
$ws.Range("B1").Value = -0.050725434693391816

$ws.Range("A2").Value = 0.022265225487358813
$ws.Range("B2").Value = -0.022265225976668333

$ws.Range("A3").Value = -0.046590928993746533
$ws.Range("B3").Value = 0.046590928477969008
